# Insert two new weekly-report rows right before the current row 1082
# (Fecha=44596 / "1a (cosecha)") so the existing data block shifts down by 2
# rows, then populate the freshly inserted rows with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1082:1083").Insert()

# Row 1082: "1a (guarda)"
$ws.Range("A1082").Value = 9
$ws.Range("B1082").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1082").Value = "Metropolitana"
$ws.Range("D1082").Value = 44769
$ws.Range("E1082").Value = 13
$ws.Range("F1082").Value = 100112045
$ws.Range("G1082").Value = "Zapallo"
$ws.Range("H1082").Value = "Camote"
$ws.Range("I1082").Value = "1a (guarda)"
$ws.Range("J1082").Value = 160
$ws.Range("K1082").Value = 880
$ws.Range("L1082").Value = 900
$ws.Range("M1082").Value = 890
$ws.Range("N1082").Value = "$/kilo (volumen en unidades)"
$ws.Range("O1082").Value = "Provincia de Melipilla"
$ws.Range("P1082").Value = 890
$ws.Range("Q1082").Value = 1
$ws.Range("R1082").Value = "Hortaliza"

# Row 1083: "2a (guarda)"
$ws.Range("A1083").Value = 9
$ws.Range("B1083").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1083").Value = "Metropolitana"
$ws.Range("D1083").Value = 44769
$ws.Range("E1083").Value = 13
$ws.Range("F1083").Value = 100112045
$ws.Range("G1083").Value = "Zapallo"
$ws.Range("H1083").Value = "Camote"
$ws.Range("I1083").Value = "2a (guarda)"
$ws.Range("J1083").Value = 70
$ws.Range("K1083").Value = 700
$ws.Range("L1083").Value = 700
$ws.Range("M1083").Value = 700
$ws.Range("N1083").Value = "$/kilo (volumen en unidades)"
$ws.Range("O1083").Value = "Provincia de Melipilla"
$ws.Range("P1083").Value = 700
$ws.Range("Q1083").Value = 1
$ws.Range("R1083").Value = "Hortaliza"

# Make sure the date cells keep the date number format used by the rest of
# column D (style index "2" in the original workbook).
$ws.Range("D1082:D1083").NumberFormat = $ws.Range("D1084").NumberFormat()
